# Apply the "alternate waterfall" regression update:
#  - shift turbine-size upsizing effect into AEP (new predictor row added)
#  - refresh the five summary coefficients (B1:B5)
#  - refresh the regression table (rows 8-16) and append a new trailing row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Top summary block (B1:B5) -------------------------------------------
$ws.Range("B1").Value = 0.7789394126229221
$ws.Range("B2").Value = 0.7356884281361025
$ws.Range("B3").Value = -0.1394830801295063
$ws.Range("B4").Value = 0.05247784092368975
$ws.Range("B5").Value = 0.09215562090148766

# --- Insert a new predictor row for "Turbine MW (Max)" --------------------
# Existing row 10 ("Distance From Shore Auto (km)") and everything below it
# shifts down by one to make room for the new predictor right after
# "Water Depth Max (m)" (row 9).
$ws.Rows.Item(10).Insert()

# New row 10: Turbine MW (Max)
$ws.Range("A10").Value = "Turbine MW (Max)"
$ws.Range("B10").Value = -0.02540563173149829
$ws.Range("C10").Value = 0.1877734100736426
$ws.Range("D10").Value = 2.020394491472369

# --- Refresh the rest of the regression table (now shifted to 11-17) -----
$ws.Range("B8").Value = 9.838244731852322
$ws.Range("C8").Value = [double]"4.04585454425828E-28"
$ws.Range("D8").Value = 293.4640298004076

$ws.Range("B9").Value = 0.002256454357066718
$ws.Range("C9").Value = 0.5141655967132666
$ws.Range("D9").Value = 3.113244435205369

$ws.Range("A11").Value = "Distance From Shore Auto (km)"
$ws.Range("B11").Value = 0.001008263652520233
$ws.Range("C11").Value = 0.3631727090827571
$ws.Range("D11").Value = 1.769866376082421

$ws.Range("A12").Value = "Germany"
$ws.Range("B12").Value = -0.09671777836262829
$ws.Range("C12").Value = 0.2103138239935414
$ws.Range("D12").Value = 2.385084263250463

$ws.Range("A13").Value = "China"
$ws.Range("B13").Value = -0.5751835369196083
$ws.Range("C13").Value = [double]"2.526736984961288E-06"
$ws.Range("D13").Value = 4.596016041490434

$ws.Range("A14").Value = "Belgium"
$ws.Range("B14").Value = -0.1914250922273117
$ws.Range("C14").Value = 0.06935966836624897
$ws.Range("D14").Value = 1.290387114164369

$ws.Range("A15").Value = "Netherlands"
$ws.Range("B15").Value = -0.186113242206102
$ws.Range("C15").Value = 0.1746956213269577
$ws.Range("D15").Value = 1.152312062370935

$ws.Range("A16").Value = "Denmark"
$ws.Range("B16").Value = -0.4300437592982931
$ws.Range("C16").Value = 0.03227215930131837
$ws.Range("D16").Value = 1.221395920549448

# --- New trailing row: log Cumulative Capacity -----------------------------
$ws.Range("A17").Value = "log Cumulative Capacity"
$ws.Range("B17").Value = -0.1394830801295063
$ws.Range("C17").Value = 0.01077644068359378
$ws.Range("D17").Value = 2.062510387941916
